$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45 (shifts existing rows 45..92 down to 46..93)
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new data record
$ws.Cells.Item(45, 1).Value  = 7
$ws.Cells.Item(45, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(45, 3).Value  = "Ñuble"
$ws.Cells.Item(45, 4).Value  = 44797
$ws.Cells.Item(45, 5).Value  = 16
$ws.Cells.Item(45, 6).Value  = 100112031
$ws.Cells.Item(45, 7).Value  = "Poroto verde"
$ws.Cells.Item(45, 8).Value  = "Magnum"
$ws.Cells.Item(45, 9).Value  = "Primera"
$ws.Cells.Item(45, 10).Value = 60
$ws.Cells.Item(45, 11).Value = 35000
$ws.Cells.Item(45, 12).Value = 35000
$ws.Cells.Item(45, 13).Value = 35000
$ws.Cells.Item(45, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(45, 15).Value = "Perú"
$ws.Cells.Item(45, 16).Value = 1400
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
